# Weekly crime data update for 79th Precinct (Volume 30, Number 18)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text runs; only replace the changed substring,
#     leaving per-run formatting of the rest of the string untouched) ---
$c = $ws.Range("A8").Characters(21, 2)
$c.Text = "18"

$c = $ws.Range("C9").Characters(27, 9)
$c.Text = "5/1/2023"
$c = $ws.Range("C9").Characters(46, 9)
$c.Text = "5/7/2023"

# --- Table data updates (rows 14-30) ---
$ws.Range("N14").Value = -90.47619047619
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 150
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 9.090909090909
$ws.Range("L15").Value = 300
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -55.555555555555
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("I16").Value = 68
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = 15.254237288135
$ws.Range("L16").Value = -15
$ws.Range("M16").Value = -41.379310344827
$ws.Range("N16").Value = -90.18759018759
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 62.5
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = 5.555555555555
$ws.Range("I17").Value = 127
$ws.Range("J17").Value = 139
$ws.Range("K17").Value = -8.633093525179
$ws.Range("L17").Value = -0.78125
$ws.Range("M17").Value = -13.013698630137
$ws.Range("N17").Value = -66.75392670157
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 72
$ws.Range("J18").Value = 73
$ws.Range("K18").Value = -1.369863013698
$ws.Range("L18").Value = 22.033898305084
$ws.Range("M18").Value = -48.936170212766
$ws.Range("N18").Value = -82.222222222222
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 60.869565217391
$ws.Range("I19").Value = 148
$ws.Range("J19").Value = 117
$ws.Range("K19").Value = 26.495726495726
$ws.Range("L19").Value = 29.824561403508
$ws.Range("M19").Value = 9.629629629629
$ws.Range("N19").Value = -43.076923076923
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "'0"
$ws.Range("E20").Value = "'***.*"
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 46
$ws.Range("K20").Value = 39.393939393939
$ws.Range("L20").Value = 43.75
$ws.Range("M20").Value = 17.948717948717
$ws.Range("N20").Value = -77.777777777777
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 28.571428571428
$ws.Range("F21").Value = 125
$ws.Range("H21").Value = 35.869565217391
$ws.Range("I21").Value = 475
$ws.Range("J21").Value = 434
$ws.Range("K21").Value = 9.447004608294
$ws.Range("L21").Value = 13.365155131264
$ws.Range("M21").Value = -19.763513513513
$ws.Range("N21").Value = -76.190476190476
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 60
$ws.Range("F23").Value = 21
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = 16.666666666666
$ws.Range("I23").Value = 83
$ws.Range("J23").Value = 67
$ws.Range("K23").Value = 23.880597014925
$ws.Range("L23").Value = -1.190476190476
$ws.Range("M23").Value = 29.6875
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 102
$ws.Range("G24").Value = 121
$ws.Range("H24").Value = -15.702479338843
$ws.Range("I24").Value = 527
$ws.Range("J24").Value = 525
$ws.Range("K24").Value = 0.380952380952
$ws.Range("L24").Value = 36.883116883116
$ws.Range("M24").Value = 70
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = -18
$ws.Range("I25").Value = 219
$ws.Range("J25").Value = 226
$ws.Range("K25").Value = -3.097345132743
$ws.Range("L25").Value = 54.225352112676
$ws.Range("M25").Value = -23.426573426573
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 18
$ws.Range("J26").Value = 17
$ws.Range("K26").Value = 5.882352941176
$ws.Range("L26").Value = 157.142857142857
$ws.Range("D27").Value = 1
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = -13.333333333333
$ws.Range("L27").Value = -40.90909090909
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = -60
$ws.Range("M28").Value = -71.428571428571
$ws.Range("N28").Value = -94.339622641509
$ws.Range("C29").Value = 2
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = -28.571428571428
$ws.Range("L29").Value = -61.538461538461
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -94.252873563218
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "'***.*"
